$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.807.13"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +2.56%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'1.876.15"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +2.17%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").Value = "'1.005"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  +0.38%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'325.07"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.37%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("E6").Value = "'  +0.27%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = "'0.4617"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -0.36%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'0.3869"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -0.20%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'0.07872"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +0.06%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'0.9863"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +2.41%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'21.88"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -0.29%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'1.900.50"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +4.66%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'7.003"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +1.09%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'5.708"
$ws.Range("D14").Style = "Normal"

$ws.Range("D15").Value = "'0.06982"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +2.33%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("E16").Value = "'  +0.21%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'1.004"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +0.24%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("E18").Value = "'  +0.98%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("E19").Value = "'  +0.47%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "'1.003"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +0.23%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'28.827.29"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +2.56%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("E22").Value = "'  -0.90%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'11.08"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +0.69%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'2.103"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +0.28%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'2.117.34"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +3.69%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'152.94"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -1.11%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = "'19.32"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +0.67%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("D28").Value = "'5.846"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +2.79%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("D29").Value = "'1.995"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +1.47%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").Value = "'118.99"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +0.77%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = "'0.09347"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +1.24%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("D32").Value = "'0.9214"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -1.67%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").Value = "'5.307"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +0.62%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("D34").Value = "'1.340"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +1.33%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("D35").Value = "'3.323"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +0.55%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("D36").Value = "'0.05795"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -1.46%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("D37").Value = "'1.150"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +0.34%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("D38").Value = "'0.02068"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -2.82%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("E39").Value = "'  -1.36%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = "'0.5634"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +0.53%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = "'0.1784"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +1.13%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("E42").Value = "'  -1.20%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'0.07218"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -0.67%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'11.77"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +0.99%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'0.5302"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +0.46%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'2.145"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +0.65%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = "'1.123"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -1.87%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'1.840"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +0.48%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value = "'113.44"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +0.79%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Value = "'2.416"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +3.77%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("E51").Value = "'  +0.26%  "
$ws.Range("E51").Style = "Normal"
